$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The CardCode (column A) and ItemCode (column B) values for every
# transaction row (rows 2-21) get an underscore prefix added, e.g.
# "CUST001" -> "_CUST001" and "ITEM0001" -> "_ITEM0001".
# Quantity (column C) and the header row (row 1) are left untouched.

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cardCode = $ws.Cells.Item($r, 1).Value2
    if ($cardCode -ne $null -and -not $cardCode.ToString().StartsWith("_")) {
        $ws.Cells.Item($r, 1).Value = "_" + $cardCode
    }

    $itemCode = $ws.Cells.Item($r, 2).Value2
    if ($itemCode -ne $null -and -not $itemCode.ToString().StartsWith("_")) {
        $ws.Cells.Item($r, 2).Value = "_" + $itemCode
    }
}

$ws.Range("B22").Select() | Out-Null
